# Generate Report for Archive
# Update the localization "Status" value from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview summary columns for
# zh-cn / de-de, and the per-language "Status" table column), then let
# the Status columns auto-size to the new (shorter) text, matching what
# Excel does automatically when the cell content driving a best-fit
# column width changes.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: per-language status summary, columns E (zh-cn) and F (de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn / de-de detail sheets: "Status" column (C)
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the columns whose width was driven by the Status text so the
# column width reflects the new, shorter value.
$overview.Columns.Item(5).EntireColumn.AutoFit()
$overview.Columns.Item(6).EntireColumn.AutoFit()
$zhcn.Columns.Item(3).EntireColumn.AutoFit()
$dede.Columns.Item(3).EntireColumn.AutoFit()
